$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Timestamp header update
$ws.Range("A1").Value = "Datos actualizados a 12 de Julio de 2020 a las 20:42"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 3390374
$ws.Range("C4").Value = 34728
$ws.Range("D4").Value = 1506094
$ws.Range("E4").Value = 1746629
$ws.Range("G4").Value = 249
$ws.Range("H4").Value = 137651

# India (row 6)
$ws.Range("B6").Value = 879447
$ws.Range("C6").Value = 29089
$ws.Range("D6").Value = 554370
$ws.Range("E6").Value = 301890
$ws.Range("G6").Value = 500
$ws.Range("H6").Value = 23187

# Chile (row 9)
$ws.Range("B9").Value = 315041
$ws.Range("C9").Value = 3012
$ws.Range("D9").Value = 283902
$ws.Range("E9").Value = 24160
$ws.Range("G9").Value = 98
$ws.Range("H9").Value = 6979

# Alemania (row 19)
$ws.Range("D19").Value = 184600
$ws.Range("E19").Value = 6180

# Barein (row 51)
$ws.Range("E51").Value = 4534
$ws.Range("G51").Value = 4
$ws.Range("H51").Value = 108

# Row 61/62: Austria & Argelia swap positions (Argelia overtakes Austria)
$ws.Range("A61").Value = "Argelia"
$ws.Range("B61").Value = 19195
$ws.Range("C61").Value = 483
$ws.Range("D61").Value = 13743
$ws.Range("E61").Value = 4441
$ws.Range("G61").Value = 7
$ws.Range("H61").Value = 1011

$ws.Range("A62").Value = "Austria"
$ws.Range("B62").Value = 18897
$ws.Range("C62").Value = 114
$ws.Range("D62").Value = 16952
$ws.Range("E62").Value = 1237
$ws.Range("G62").Value = 2
$ws.Range("H62").Value = 708

# Uzbekistan (row 70)
$ws.Range("E70").Value = 5072
$ws.Range("G70").Value = 3
$ws.Range("H70").Value = 60

# Estado de Palestina (row 90)
$ws.Range("E90").Value = 5252
$ws.Range("G90").Value = 3
$ws.Range("H90").Value = 36

# Row 108/109: Mayotte & Maldivas swap positions (Maldivas overtakes Mayotte)
$ws.Range("A108").Value = "Maldivas"
$ws.Range("B108").Value = 2731
$ws.Range("C108").Value = 67
$ws.Range("D108").Value = 2284
$ws.Range("E108").Value = 434
$ws.Range("H108").Value = 13

$ws.Range("A109").Value = "Mayotte"
$ws.Range("B109").Value = 2711
$ws.Range("D109").Value = 2480
$ws.Range("E109").Value = 194
$ws.Range("H109").Value = 37

# Sri Lanka (row 110)
$ws.Range("B110").Value = 2617
$ws.Range("C110").Value = 106
$ws.Range("E110").Value = 625

# Row 113/114: Libano & Malaui swap positions (Malaui overtakes Libano)
$ws.Range("A113").Value = "Malaui"
$ws.Range("B113").Value = 2364
$ws.Range("C113").Value = 103
$ws.Range("D113").Value = 557
$ws.Range("E113").Value = 1769
$ws.Range("G113").Value = 5
$ws.Range("H113").Value = 38

$ws.Range("A114").Value = "Libano"
$ws.Range("B114").Value = 2334
$ws.Range("C114").Value = 166
$ws.Range("D114").Value = 1420
$ws.Range("E114").Value = 878
$ws.Range("H114").Value = 36

# Row 128/129: Libia & Yemen swap positions (Yemen overtakes Libia)
$ws.Range("A128").Value = "Yemen"
$ws.Range("B128").Value = 1465
$ws.Range("C128").Value = 76
$ws.Range("D128").Value = 659
$ws.Range("E128").Value = 389
$ws.Range("G128").Value = 52
$ws.Range("H128").Value = 417

$ws.Range("A129").Value = "Libia"
$ws.Range("B129").Value = 1433
$ws.Range("C129").Value = 44
$ws.Range("D129").Value = 341
$ws.Range("E129").Value = 1053
$ws.Range("G129").Value = 1
$ws.Range("H129").Value = 39

# Row 131
$ws.Range("B131").Value = 1351
$ws.Range("C131").Value = 40
$ws.Range("D131").Value = 668
$ws.Range("E131").Value = 663
$ws.Range("G131").Value = 2
$ws.Range("H131").Value = 20
